$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamp to stamp on every data row (O column = "timestamp")
$newTimestamp = "2022-12-26 20:52:20"

# Determine the last used row from column A (id column), falling back to
# the known data extent if detection ever comes back empty.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 547 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# Update the two product titles that went out of stock online
$ws.Cells.Item(184, 13).Value = "Prix Garantie Champignons - Online kein Bestand 2.00 Schweizer Franken"
$ws.Cells.Item(334, 13).Value = "Naturaplan Bio " + [char]0x00C4 + "pfel rote Sorte 750g - Online kein Bestand 4.95 Schweizer Franken"
